$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy column D's formatting onto the new column E (rows 4-43) ---
$ws.Range("D4:D43").Copy()
$ws.Range("E4:E43").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Header row: row height grows to fit the longer footnote ---
$ws.Rows.Item(1).RowHeight = 67.5

# --- Column widths: A,B,C all become a uniform 44 ---
$ws.Columns.Item(1).ColumnWidth = 43.17
$ws.Columns.Item(2).ColumnWidth = 43.17
$ws.Columns.Item(3).ColumnWidth = 43.17

# --- Row 4: survey year header gains a 2023 column ---
$ws.Range("E4").Value = 2023

# --- Row 5: Total row gains 2023 total ---
$ws.Range("E5").Value = 6.2

# --- Row 6: Urbanisation group header (no data) ---
$ws.Range("E6").ClearContents()

# --- Row 7-8: City / Village ---
$ws.Range("E7").Value = 7.4
$ws.Range("E8").Value = 5.6

# --- Row 9: By territory group header (no data) ---
$ws.Range("E9").ClearContents()

# --- Rows 10-18: oblasts ---
$ws.Range("E10").Value = 4.3
$ws.Range("E11").Value = 7.1
$ws.Range("E12").Value = 2.5
$ws.Range("E13").Value = 2.9
$ws.Range("E14").Value = 3.4
$ws.Range("E15").Value = 1.9
$ws.Range("E16").Value = 9.3
$ws.Range("E17").Value = 7.1
$ws.Range("E18").Value = 14.9

# --- Row 19: Age group header (no data) ---
$ws.Range("E19").ClearContents()

# --- Rows 20-28: age bands ---
$ws.Range("E20").Value = 5.3
$ws.Range("E21").Value = 3.5
$ws.Range("E22").Value = 10
$ws.Range("E23").Value = 5.3
$ws.Range("E24").Value = 5.5
$ws.Range("E25").Value = 7.7
$ws.Range("E26").Value = 6.8
$ws.Range("E27").Value = 5.8
$ws.Range("E28").Value = 7

# --- Row 29: Education group header (no data) ---
$ws.Range("E29").ClearContents()

# --- Row 30: Preschool/none/primary row — footnote marker instead of a number ---
$ws.Range("E30").Value = "(18,7)"

# --- Rows 31-34: remaining education levels ---
$ws.Range("E31").Value = 7.5
$ws.Range("E32").Value = 6.1
$ws.Range("E33").Value = 4.9
$ws.Range("E34").Value = 6.4

# --- Row 35: Functional difficulties group header (no data) ---
$ws.Range("E35").ClearContents()

# --- Rows 36-37: Yes / No ---
$ws.Range("E36").Value = 32.3
$ws.Range("E37").Value = 6.1

# --- Row 38: Wealth quintile group header (no data) ---
$ws.Range("E38").ClearContents()

# --- Rows 39-43: wealth quintiles ---
$ws.Range("E39").Value = 6.7
$ws.Range("E40").Value = 5.7
$ws.Range("E41").Value = 5
$ws.Range("E42").Value = 6.4
$ws.Range("E43").Value = 7.1

# --- Footnote row (44): update text in all three languages to mention 2023 too ---
$ws.Range("A44").Value = " Көп көрсөткүчтүү кластердик изилдөөнүн маалыматтары боюнча, 2018-ж., 2023-ж."
$ws.Range("B44").Value = "По данным кластерного обследования по многим показателям, 2018г., 2023г."
$ws.Range("C44").Value = "According to Multiple Indicator Cluster Survey, 2018, 2023."
